{"js": "async (context) => {\n  const replacements = [\n    [\"23\u00d726=598\", \"46\u00d749=2254\"],\n    [\"39\u00d784=3276\", \"44\u00d718=792\"],\n    [\"21\u00d786=1806\", \"57\u00d757=3249\"],\n    [\"79\u00d789=7031\", \"20\u00d718=360\"],\n    [\"46\u00d778=3588\", \"22\u00d724=528\"],\n    [\"31\u00d794=2914\", \"96\u00d726=2496\"],\n    [\"69\u00d786=5934\", \"46\u00d737=1702\"],\n    [\"32\u00d753=1696\", \"28\u00d795=2660\"],\n    [\"55\u00d794=5170\", \"29\u00d772=2088\"],\n    [\"90\u00d782=7380\", \"61\u00d754=3294\"],\n    [\"55\u00d770=3850\", \"89\u00d769=6141\"],\n    [\"46\u00d750=2300\", \"48\u00d748=2304\"],\n    [\"54\u00d720=1080\", \"56\u00d731=1736\"],\n    [\"98\u00d721=2058\", \"19\u00d780=1520\"],\n    [\"14\u00d716=224\", \"43\u00d758=2494\"],\n    [\"70\u00d733=2310\", \"42\u00d737=1554\"],\n    [\"45\u00d795=4275\", \"85\u00d721=1785\"],\n    [\"41\u00d789=3649\", \"30\u00d742=1260\"],\n    [\"22\u00d782=1804\", \"37\u00d777=2849\"],\n    [\"18\u00d738=684\", \"18\u00d793=1674\"],\n    [\"46\u00d794=4324\", \"82\u00d761=5002\"],\n    [\"87\u00d787=7569\", \"34\u00d755=1870\"],\n    [\"34\u00d792=3128\", \"82\u00d760=4920\"],\n    [\"93\u00d724=2232\", \"55\u00d747=2585\"],\n    [\"12\u00d741=492\", \"73\u00d761=4453\"],\n  ];\n\n  const body = context.document.body;\n\n  for (const [oldText, newText] of replacements) {\n    const searchResults = body.search(oldText, { matchCase: true, matchWholeWord: false });\n    searchResults.load(\"items\");\n    await context.sync();\n\n    for (const range of searchResults.items) {\n      range.insertText(newText, Word.InsertLocation.replace);\n    }\n    await context.sync();\n  }\n};\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"23\u00d726=598\", \"46\u00d749=2254\"),\n    @(\"39\u00d784=3276\", \"44\u00d718=792\"),\n    @(\"21\u00d786=1806\", \"57\u00d757=3249\"),\n    @(\"79\u00d789=7031\", \"20\u00d718=360\"),\n    @(\"46\u00d778=3588\", \"22\u00d724=528\"),\n    @(\"31\u00d794=2914\", \"96\u00d726=2496\"),\n    @(\"69\u00d786=5934\", \"46\u00d737=1702\"),\n    @(\"32\u00d753=1696\", \"28\u00d795=2660\"),\n    @(\"55\u00d794=5170\", \"29\u00d772=2088\"),\n    @(\"90\u00d782=7380\", \"61\u00d754=3294\"),\n    @(\"55\u00d770=3850\", \"89\u00d769=6141\"),\n    @(\"46\u00d750=2300\", \"48\u00d748=2304\"),\n    @(\"54\u00d720=1080\", \"56\u00d731=1736\"),\n    @(\"98\u00d721=2058\", \"19\u00d780=1520\"),\n    @(\"14\u00d716=224\", \"43\u00d758=2494\"),\n    @(\"70\u00d733=2310\", \"42\u00d737=1554\"),\n    @(\"45\u00d795=4275\", \"85\u00d721=1785\"),\n    @(\"41\u00d789=3649\", \"30\u00d742=1260\"),\n    @(\"22\u00d782=1804\", \"37\u00d777=2849\"),\n    @(\"18\u00d738=684\", \"18\u00d793=1674\"),\n    @(\"46\u00d794=4324\", \"82\u00d761=5002\"),\n    @(\"87\u00d787=7569\", \"34\u00d755=1870\"),\n    @(\"34\u00d792=3128\", \"82\u00d760=4920\"),\n    @(\"93\u00d724=2232\", \"55\u00d747=2585\"),\n    @(\"12\u00d741=492\", \"73\u00d761=4453\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute(\n        $oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2\n    ) | Out-Null\n}\n"}
